$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4487.08
$ws.Range("J17").Value = 4487.08
$ws.Range("L17").Value = 13461.24
$ws.Range("N17").Value = -13797.24
$ws.Range("H19").Value = 737.7
$ws.Range("I19").Value = 749.5
$ws.Range("J19").Value = 729.8333
$ws.Range("K19").Value = 749.5
$ws.Range("L19").Value = 729.8333
$ws.Range("M19").Value = -574.5
$ws.Range("N19").Value = -1079.8333
$ws.Range("H28").Value = 1076.6
$ws.Range("I28").Value = 1076.6
$ws.Range("K28").Value = 1076.6
$ws.Range("M28").Value = -591.5999999999999
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H43").Value = 2411.8572
$ws.Range("I43").Value = 1628
$ws.Range("K43").Value = 1628
$ws.Range("M43").Value = -1559
$ws.Range("H51").Value = 8000.5
$ws.Range("I51").Value = 8000.5
$ws.Range("K51").Value = 8000.5
$ws.Range("M51").Value = -7516.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7989.875
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7989.875
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -9237.875
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7989.875
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 39949.375
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -46189.375
$ws.Range("H98").Value = 3018.5715
$ws.Range("I98").Value = 917.44446
$ws.Range("J98").Value = 6800.6
$ws.Range("K98").Value = 917.44446
$ws.Range("L98").Value = 6800.6
$ws.Range("M98").Value = 580.55554
$ws.Range("N98").Value = -9796.6
$ws.Range("H118").Value = 3522.5
$ws.Range("I118").Value = 3522.5
$ws.Range("K118").Value = 10567.5
$ws.Range("M118").Value = -8910.5
$ws.Range("H122").Value = 3018.5715
$ws.Range("I122").Value = 917.44446
$ws.Range("J122").Value = 6800.6
$ws.Range("K122").Value = 2752.33338
$ws.Range("L122").Value = 20401.8
$ws.Range("M122").Value = -302.33338
$ws.Range("N122").Value = -25301.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6554.9546
$ws.Range("I32").Value = 5914.7617
$ws.Range("K32").Value = 5914.7617
$ws.Range("M32").Value = -5627.7617
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1623
$ws.Range("N45").ClearContents()
$ws.Range("H102").Value = 794.8
$ws.Range("I102").Value = 794.8
$ws.Range("K102").Value = 794.8
$ws.Range("M102").Value = 827.2
$ws.Range("H122").Value = 1668.8125
$ws.Range("I122").Value = 956.875
$ws.Range("J122").Value = 2380.75
$ws.Range("K122").Value = 2870.625
$ws.Range("L122").Value = 7142.25
$ws.Range("M122").Value = -420.625
$ws.Range("N122").Value = -12042.25
$ws.Range("H132").Value = 1949.6
$ws.Range("J132").Value = 2751
$ws.Range("L132").Value = 8253
$ws.Range("N132").Value = -13313

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5030.231
$ws.Range("I94").Value = 6289.3
$ws.Range("K94").Value = 6289.3
$ws.Range("M94").Value = -5838.3
$ws.Range("H99").Value = 1708.2273
$ws.Range("I99").Value = 1119.5
$ws.Range("J99").Value = 2738.5
$ws.Range("K99").Value = 1119.5
$ws.Range("L99").Value = 2738.5
$ws.Range("M99").Value = 378.5
$ws.Range("N99").Value = -5734.5
$ws.Range("H114").Value = 44444
$ws.Range("J114").Value = 44444
$ws.Range("L114").Value = 44444
$ws.Range("N114").Value = -53122

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2863.6667
$ws.Range("I58").Value = 2837
$ws.Range("J58").Value = 2997
$ws.Range("K58").Value = 2837
$ws.Range("L58").Value = 2997
$ws.Range("M58").Value = -2634
$ws.Range("N58").Value = -3403
$ws.Range("H86").Value = 11939
$ws.Range("I86").Value = 12421.75
$ws.Range("K86").Value = 12421.75
$ws.Range("M86").Value = -11298.75
$ws.Range("H89").Value = 11939
$ws.Range("I89").Value = 12421.75
$ws.Range("K89").Value = 62108.75
$ws.Range("M89").Value = -56492.75
$ws.Range("H94").Value = 163747.86
$ws.Range("I94").Value = 552512
$ws.Range("K94").Value = 552512
$ws.Range("M94").Value = -552061
$ws.Range("H122").Value = 1451.8462
$ws.Range("I122").Value = 1037.5714
$ws.Range("K122").Value = 3112.7142
$ws.Range("M122").Value = -662.7142000000003
$ws.Range("H134").Value = 3107.3076
$ws.Range("I134").Value = 3212.4443
$ws.Range("J134").Value = 2870.75
$ws.Range("K134").Value = 9637.332900000001
$ws.Range("L134").Value = 8612.25
$ws.Range("M134").Value = -7102.332900000001
$ws.Range("N134").Value = -13682.25
$ws.Range("H136").Value = 2863.6667
$ws.Range("I136").Value = 2837
$ws.Range("J136").Value = 2997
$ws.Range("K136").Value = 8511
$ws.Range("L136").Value = 8991
$ws.Range("M136").Value = -5961
$ws.Range("N136").Value = -14091

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13333520
$ws.Range("I4").Value = 13750210
$ws.Range("K4").Value = 41250630
$ws.Range("M4").Value = -41250518
$ws.Range("H11").Value = 53609412
$ws.Range("I11").Value = 53609412
$ws.Range("K11").Value = 160828236
$ws.Range("M11").Value = -160828096
$ws.Range("H39").Value = 7059.6
$ws.Range("J39").Value = 6824.5
$ws.Range("L39").Value = 20473.5
$ws.Range("N39").Value = -21061.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1312.75
$ws.Range("I102").Value = 1215.4
$ws.Range("J102").Value = 1799.5
$ws.Range("K102").Value = 1215.4
$ws.Range("L102").Value = 1799.5
$ws.Range("M102").Value = 406.5999999999999
$ws.Range("N102").Value = -5043.5
$ws.Range("H113").Value = 1684.4286
$ws.Range("I113").Value = 1631.8334
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1631.8334
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 538.1666
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 6014.273
$ws.Range("I122").Value = 4807.4287
$ws.Range("J122").Value = 8126.25
$ws.Range("K122").Value = 14422.2861
$ws.Range("L122").Value = 24378.75
$ws.Range("M122").Value = -11972.2861
$ws.Range("N122").Value = -29278.75
$ws.Range("H126").Value = 3620.625
$ws.Range("I126").Value = 1817
$ws.Range("J126").Value = 4702.8
$ws.Range("K126").Value = 5451
$ws.Range("L126").Value = 14108.4
$ws.Range("M126").Value = -2981
$ws.Range("N126").Value = -19048.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 132500
$ws.Range("I2").Value = 132500
$ws.Range("K2").Value = 132500
$ws.Range("M2").Value = -132388
$ws.Range("H7").Value = 6505.857
$ws.Range("I7").Value = 2319.889
$ws.Range("K7").Value = 2319.889
$ws.Range("M7").Value = -2207.889
$ws.Range("H16").Value = 2836.375
$ws.Range("I16").Value = 3115.6667
$ws.Range("J16").Value = 1998.5
$ws.Range("K16").Value = 3115.6667
$ws.Range("L16").Value = 1998.5
$ws.Range("M16").Value = -2945.6667
$ws.Range("N16").Value = -2338.5
$ws.Range("H22").Value = 862.25
$ws.Range("I22").Value = 862.25
$ws.Range("K22").Value = 862.25
$ws.Range("M22").Value = -567.25
$ws.Range("H27").Value = 862.25
$ws.Range("I27").Value = 862.25
$ws.Range("K27").Value = 862.25
$ws.Range("M27").Value = -755.25
$ws.Range("H40").Value = 3116.25
$ws.Range("I40").Value = 2323.125
$ws.Range("J40").Value = 3909.375
$ws.Range("K40").Value = 2323.125
$ws.Range("L40").Value = 3909.375
$ws.Range("M40").Value = -2187.125
$ws.Range("N40").Value = -4181.375
$ws.Range("H61").Value = 4698.7
$ws.Range("J61").Value = 5074.5
$ws.Range("L61").Value = 5074.5
$ws.Range("N61").Value = -5478.5
$ws.Range("H113").Value = 4698.7
$ws.Range("J113").Value = 5074.5
$ws.Range("L113").Value = 5074.5
$ws.Range("N113").Value = -9414.5
$ws.Range("H126").Value = 6505.857
$ws.Range("I126").Value = 2319.889
$ws.Range("K126").Value = 6959.667
$ws.Range("M126").Value = -4489.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10749.75
$ws.Range("I62").Value = 6501
$ws.Range("J62").Value = 14998.5
$ws.Range("K62").Value = 6501
$ws.Range("L62").Value = 14998.5
$ws.Range("M62").Value = -5877
$ws.Range("N62").Value = -16246.5
$ws.Range("H65").Value = 10749.75
$ws.Range("I65").Value = 6501
$ws.Range("J65").Value = 14998.5
$ws.Range("K65").Value = 32505
$ws.Range("L65").Value = 74992.5
$ws.Range("M65").Value = -29385
$ws.Range("N65").Value = -81232.5
$ws.Range("H100").Value = 50005000
$ws.Range("J100").Value = 10000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -21082
$ws.Range("H107").Value = 681
$ws.Range("I107").Value = 86
$ws.Range("K107").Value = 258
$ws.Range("M107").Value = 1662
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 402.33334
$ws.Range("I122").Value = 402.33334
$ws.Range("K122").Value = 1207.00002
$ws.Range("M122").Value = 1242.99998
$ws.Range("H124").Value = 24998.5
$ws.Range("J124").Value = 24998.5
$ws.Range("L124").Value = 24998.5
$ws.Range("N124").Value = -34818.5
$ws.Range("H126").Value = 7134.5713
$ws.Range("I126").Value = 6154.222
$ws.Range("K126").Value = 18462.666
$ws.Range("M126").Value = -15992.666
